$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "total score" formula to column E (E2:E146): 0.6*C + 0.4*D
$ws.Range("E2:E146").Formula = "=0.6*C2+0.4*D2"

# Restore the selection recorded by Excel when the author scrolled
# down to the bottom of the sheet after filling in the formulas
$ws.Range("G145").Select()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
